$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BrokenMedia")

# Update browserId (column A) and parentPage (column B) values per corrected mapping.
$ws.Range("A2").Value = "2-1"
$ws.Range("B2").Value = "http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/"

$ws.Range("A3").Value = "2-2"
$ws.Range("B3").Value = "http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/"

$ws.Range("A4").Value = "2-2"

$ws.Range("A5").Value = "2-1"
